$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + "30.523.21"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.60%  "
$ws.Range("D3").Value = "'" + "1.853.27"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.28%  "
$ws.Range("D4").Value = "'" + "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'" + "233.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.21%  "
$ws.Range("D6").Value = "'" + "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("D7").Value = "'" + "0.4747"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.25%  "
$ws.Range("D8").Value = "'" + "0.2745"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.60%  "
$ws.Range("D9").Value = "'" + "0.06324"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.54%  "
$ws.Range("D10").Value = "'" + "17.65"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +11.42%  "
$ws.Range("D11").Value = "'" + "1.841.99"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.67%  "
$ws.Range("D12").Value = "'" + "0.07450"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.68%  "
$ws.Range("D13").Value = "'" + "4.968"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.12%  "
$ws.Range("D14").Value = "'" + "84.63"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.97%  "
$ws.Range("D15").Value = "'" + "0.6257"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.96%  "
$ws.Range("D16").Value = "'" + "30.516.98"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.79%  "
$ws.Range("D17").Value = "'" + "247.09"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +10.27%  "
$ws.Range("E18").Value = "  -0.13%  "
$ws.Range("D19").Value = "'" + "12.71"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.32%  "
$ws.Range("D20").Value = "'" + "0.000007329"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.02%  "
$ws.Range("D21").Value = "'" + "1.000"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.39%  "
$ws.Range("D22").Value = "'" + "4.929"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.91%  "
$ws.Range("D23").Value = "'" + "5.912"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.83%  "
$ws.Range("D24").Value = "'" + "9.131"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.20%  "
$ws.Range("D25").Value = "'" + "164.21"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("D26").Value = "'" + "17.96"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.90%  "
$ws.Range("D27").Value = "'" + "1.871"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.37%  "
$ws.Range("D28").Value = "'" + "0.1026"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.15%  "
$ws.Range("D29").Value = "'" + "1.359"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.66%  "
$ws.Range("E30").Value = "  -0.07%  "
$ws.Range("D31").Value = "'" + "3.836"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.01%  "
$ws.Range("D32").Value = "'" + "0.04846"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.67%  "
$ws.Range("D33").Value = "'" + "1.133"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.86%  "
$ws.Range("D34").Value = "'" + "0.6979"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.30%  "
$ws.Range("D35").Value = "'" + "2.701"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.41%  "
$ws.Range("D36").Value = "'" + "0.01898"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.24%  "
$ws.Range("D37").Value = "'" + "2.680"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.97%  "
$ws.Range("D38").Value = "'" + "2.005"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.05%  "
$ws.Range("D39").Value = "'" + "0.8743"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.52%  "
$ws.Range("D40").Value = "'" + "106.32"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.76%  "
$ws.Range("D41").Value = "'" + "1.001"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.49%  "
$ws.Range("D42").Value = "'" + "5.538"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.38%  "
$ws.Range("D43").Value = "'" + "0.4063"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.81%  "
$ws.Range("D44").Value = "'" + "7.192"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.29%  "
$ws.Range("D45").Value = "'" + "63.06"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.57%  "
$ws.Range("D46").Value = "'" + "0.1201"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.46%  "
$ws.Range("D47").Value = "'" + "33.65"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.35%  "
$ws.Range("D48").Value = "'" + "8.583"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.45%  "
$ws.Range("D49").Value = "'" + "0.05524"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("D50").Value = "'" + "1.353"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.77%  "
$ws.Range("D51").Value = "'" + "0.3690"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.17%  "
